$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("AmortTemplateGrid")
$ws2.AutoFilterMode = $false
$rng = $ws2.Range("A1:O125")
$rng.EntireRow.Hidden = $false

$ws1 = $wb.Worksheets.Item("AmortTemplateSectionGrid")
$ws1.Activate()
$ws1.Range("F2").Select()
Write-Host "done"
